# Correção das notas do fórum para matc65 em 2021.2
# Para cada aluno cuja nota_view (coluna J) seja 4, zera as colunas B:J
# (respostas diárias, total_views e nota_view).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $notaView = $ws.Cells.Item($r, 10).Value()
    if ($notaView -eq 4) {
        for ($c = 2; $c -le 10; $c++) {
            $ws.Cells.Item($r, $c).Value = 0
        }
    }
}
